# Roboflow Annotation Report 6/20/2025 - Good Night
# Fill in the labeling progress numbers for the 20/6/2025 row (row 41)
# and leave the selection positioned on the next (blank) row, as the
# author left the workbook after entering the day's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("E41").Value = 297
$ws.Range("F41").Value = 629
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 1012
$ws.Range("J41").Value = "Seguire trabajando en fin de semana (rafael)"

# Scroll the viewport so row 25 / column E is the top-left visible cell,
# then move the active selection down to D42.
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("D42").Select()
